{"js": "// The author re-typed the swear word \"F*CK\" into \"DAMN\" at the start of\n// Rick's \"F*CK IT MORTY, ...\" line. Word tracks the insertion point with\n// the special \"_GoBack\" bookmark, which moves from its old location\n// (end of the title paragraph) to right after the newly typed word.\n\nconst body = context.document.body;\n\n// 1) Drop the old \"_GoBack\" bookmark (sits at the end of the title line).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Replace \"F*CK\" with \"DAMN\" in the Rick line.\nconst results = body.search(\"F*CK\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  const newRange = target.insertText(\"DAMN\", \"Replace\");\n  await context.sync();\n\n  // 3) Re-insert \"_GoBack\" right after \"DAMN\" (before the trailing\n  //    \" IT MORTY, ...\" text), splitting the run exactly like Word does.\n  const endOfNewText = newRange.getRange(\"End\");\n  endOfNewText.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author re-typed the swear word \"F*CK\" into \"DAMN\" at the start of\n# Rick's \"F*CK IT MORTY, ...\" line. Word tracks the insertion point with\n# the special \"_GoBack\" bookmark, which moves from its old location\n# (end of the title paragraph) to right after the newly typed word.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark (sits at the end of the title line).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Replace \"F*CK\" with \"DAMN\" in the Rick line.\n$rng = $d.Content\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute(\"F*CK\")\n\nif ($found) {\n    $startPos = $rng.Start\n    $rng.Text = \"DAMN\"\n    $newEnd = $startPos + 4\n\n    # 3) Re-insert \"_GoBack\" right after \"DAMN\" (before the trailing\n    #    \" IT MORTY, ...\" text), splitting the run exactly like Word does.\n    $bmRange = $d.Range($newEnd, $newEnd)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
